# Excel COM-interop script implementing:
#   "Add a new way to work with the timezone & modify graphs & automatisation"
#
# Concretely, against plots/alerts_ech.xlsx (Sheet1, a time-series table with
# one header row and 11 data rows ECHL01.. / "Bilan en kWh"):
#   1. Fix the label for alert A9: "Jeudi1 Septembre 6h - 12h" -> "Jeudi1 Septembre 17h - 23h"
#   2. Append two new alert blocks (A15 and A16), each 7 columns wide
#      (-12h,-6h,-3h,<event label>,+3h,+6h,+12h), with empty numeric data cells
#      for every existing data row (the new alerts have not been measured yet).
#   3. Widen a couple of columns that were tight, and update the saved
#      selection / conditional-formatting range to match the new extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Correct the mislabeled time window for alert A9.
# ---------------------------------------------------------------------------
$ws.Range("BI1").Value = "Jeudi1 Septembre 17h - 23h"

# ---------------------------------------------------------------------------
# 2) Append alert blocks A15 (CV:DJ.. actually CV:DB) and A16 (DC:DI).
#    Existing table runs from column B to column CU (alerts A1..A14).
#    Each alert occupies 7 columns: -12h -6h -3h <event> +3h +6h +12h.
# ---------------------------------------------------------------------------
$headers = @(
    "A15-12h", "A15-6h", "A15-3h", "Dimanche4 Decembre 12h - 18h", "A15+3h", "A15+6h", "A15+12h",
    "A16-12h", "A16-6h", "A16-3h", "Lundi19 Decembre 16h - 22h", "A16+3h", "A16+6h", "A16+12h"
)
$cols = @("CV", "CW", "CX", "CY", "CZ", "DA", "DB", "DC", "DD", "DE", "DF", "DG", "DH", "DI")

# Clone the header-row style (bold, centered, bordered) from the last existing
# header cell so the new headers look consistent instead of falling back to
# the default style.
$ws.Range("CU1").Copy()
$ws.Range("CV1:DI1").PasteSpecial(-4122)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# Clone the numeric-cell style (2 decimal places) down each data row so the
# new, still-empty cells match the rest of the table.
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("CU" + $r).Copy()
    $ws.Range("CV" + $r + ":DI" + $r).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 3) Cosmetic follow-up: widen a couple of columns, refresh the stored
#    selection, and extend the conditional color-scale formatting to the
#    new range.
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 9.33
$ws.Columns.Item(12).ColumnWidth = 7

$ws.Range("CW18").Select()

$cf = $ws.Range("A1:CU11").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A1:DI11"))

Write-Output "alerts_ech edit applied"
